$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verleiherabgaben")

# Delete the three rows (in the original layout: rows 6, 8 and 10) that held
# "Wisdom of Happiness - A heart-to-heart with the Dalai Lama" (Elite Film AG),
# "A Real Pain" (The Walt Disney Company Switzerland GmbH) and
# "Lee - Die Fotografin" (Elite Film AG). Deleting from the bottom up keeps the
# remaining row numbers stable while we work.
$ws.Rows(10).Delete()
$ws.Rows(8).Delete()
$ws.Rows(6).Delete()

$ws.Range("D25").Select()
